# Update stats for 2025-08 (row 21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New source values
$ws.Range("B21").Value2 = 6207
$ws.Range("D21").Value2 = 5582477

# Recompute dependent metrics to stay consistent with the rest of the sheet
$B9 = $ws.Range("B9").Value2
$D9 = $ws.Range("D9").Value2

$B21 = $ws.Range("B21").Value2
$D21 = $ws.Range("D21").Value2

$ws.Range("E21").Value2 = $D21 / $B21
$ws.Range("F21").Value2 = ($B21 - $B9) / $B9 * 100
$ws.Range("H21").Value2 = ($D21 - $D9) / $D9 * 100
